$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The D-column values are plain text like "01/01/2023", not real dates
# (the source cells have no date number-format). Assigning such a string
# directly makes Excel auto-convert it to a date serial, so we briefly
# force a text format, assign the value, then restore the default
# "Normal" style so the resulting cell has no style override (matching
# the original, unstyled cells).

# Row 32: keep C, change D (date) and E (value)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "01/01/2023"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = 2.276633442999038

# Row 33: change C, D, E
$ws.Range("C33").Value = "Administração Pública"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "01/01/2023"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = 2.315008864560557

# Row 34: change C, D, E
$ws.Range("C34").Value = "Entidades Empresariais"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "01/01/2023"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = 40.47589794891648

# Row 35: change C, E (D already 01/01/2023)
$ws.Range("C35").Value = "Entidades Empresariais"
$ws.Range("E35").Value = 40.48422032949609

# Row 36: change C, E (D already 01/01/2023)
$ws.Range("C36").Value = "Entidades sem Fins Lucrativos"
$ws.Range("E36").Value = 22.06448033444356
